# daily auto push: 2026-02-26 14:14 UTC
# Two new rows of data for 2026/02/26 (time slots 20 and 22) are inserted
# right after the existing 2026/02/26 entries (which end at row 879),
# pushing every following row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 880:881 - everything from the old row 880
# onward (2026/12/29 ...) shifts down to 882 onward.
$ws.Rows("880:881").Insert()

# Row 880: 2026/02/26, 木, 20, 201
$c = $ws.Cells.Item(880, 1)
$c.NumberFormat = "@"        # force text so the date-like string isn't
$c.Value2 = "2026/02/26"     # auto-converted into a date serial value
$c.ClearFormats()            # drop the temporary text format again so the
                              # cell ends up with no explicit style, like its
                              # neighbours
$ws.Cells.Item(880, 2).Value2 = "木"
$ws.Cells.Item(880, 3).Value2 = 20
$ws.Cells.Item(880, 4).Value2 = 201

# Row 881: 2026/02/26, 木, 22, 201
$c2 = $ws.Cells.Item(881, 1)
$c2.NumberFormat = "@"
$c2.Value2 = "2026/02/26"
$c2.ClearFormats()
$ws.Cells.Item(881, 2).Value2 = "木"
$ws.Cells.Item(881, 3).Value2 = 22
$ws.Cells.Item(881, 4).Value2 = 201
